$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5's numeric values to 2 decimal places (custom accuracy).
$ws.Range("B5").Value = 22.6
$ws.Range("C5").Value = 17.23
$ws.Range("D5").Value = 0.71
$ws.Range("E5").Value = 48.76
$ws.Range("F5").Value = 40.68
$ws.Range("G5").Value = 17.43
$ws.Range("I5").Value = 27.09
$ws.Range("J5").Value = 12.19
$ws.Range("K5").Value = 18.07
$ws.Range("L5").Value = 20.09
$ws.Range("M5").Value = 21.02
$ws.Range("N5").Value = 5.77
$ws.Range("O5").Value = 17.52
$ws.Range("P5").Value = 25.45
$ws.Range("Q5").Value = 14.62
$ws.Range("R5").Value = 0.5
$ws.Range("S5").Value = 0.54
$ws.Range("T5").Value = 261.45
$ws.Range("V5").Value = 16.43
$ws.Range("X5").Value = 17.4
$ws.Range("AA5").Value = 14.4
$ws.Range("AB5").Value = 12.5
$ws.Range("AC5").Value = 14.91
$ws.Range("AD5").Value = 20.67
$ws.Range("AF5").Value = 61.91
$ws.Range("AG5").Value = 9.82
$ws.Range("AH5").Value = 20.22

# Remove the last data row (row 6) entirely.
$ws.Rows.Item(6).Delete()
